$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Set the value of D7 to the text "n" (stored as shared string)
$ws.Range("D7").Value = "n"

# Move the active selection to D8, matching the saved view state
$ws.Range("D8").Select()
